$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column D ("Ngày nhận hàng") and shift remaining columns left
$ws.Range("D:D").Delete()

# Update the view: zoom and selection per the author's final state
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("D16").Select()
